$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1.x rows: add "Assigned to" (col C) and "Duration" (col D) ---
$ws.Range("C5").Value = "Dau"
$ws.Range("D5").Value = "2  day"

$ws.Range("C6").Value = "Dau"
$ws.Range("D6").Value = "2  day"

$ws.Range("C7").Value = "Dau"
$ws.Range("D7").Value = "1 day"

$ws.Range("C8").Value = "Tri"
$ws.Range("D8").Value = "1 day"

# --- Section 2.1.x rows (KSA / Init S and T array, Permutes S): ---
# add assignee (col C) and Start/End dates (col D/E), formatted as "d-mmm"
$ws.Range("C11").Value = "Truong"
$ws.Range("D11").Value = 44293
$ws.Range("D11").NumberFormat = "d-mmm"
$ws.Range("E11").Value = 44293
$ws.Range("E11").NumberFormat = "d-mmm"

$ws.Range("C12").Value = "Truong"
$ws.Range("D12").Value = 44293
$ws.Range("D12").NumberFormat = "d-mmm"
$ws.Range("E12").Value = 44293
$ws.Range("E12").NumberFormat = "d-mmm"

# --- Section 2.2.1 (Generate keystream) ---
$ws.Range("C15").Value = "Tri"
$ws.Range("D15").Value = "1 day"

# --- Section 2.3.1 (Convert to binary) ---
$ws.Range("C18").Value = "Dau"
$ws.Range("D18").Value = "2  day"

# --- Section 2.3.2 (row 19, "Swap") got removed entirely ---
$ws.Range("A19").ClearContents()
$ws.Range("B19").ClearContents()

# --- Section 2.4 (Excute) keeps its label, gains assignee + duration ---
$ws.Range("C21").Value = "Tri"
$ws.Range("D21").Value = "1 day"

# Final selected cell, matching the last cell the author was working on
[void]$ws.Range("C15").Select()
